# "Updates based on Dave's comments" -- fix two text issues on slide 3 of
# the high-level diagram deck, and (best effort) refresh the auto date
# field on the Handout Master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Handout Master date placeholder: "11/19/20" -> "1/12/21".
#    This text lives inside an auto-updating
#    <a:fld type="datetimeFigureOut"> field, which PowerPoint
#    recomputes from the system clock on save rather than accepting a
#    literal string from automation. We still call the documented
#    Headers/Footers API for it; on hosts where that doesn't stick this
#    is a harmless no-op and every other edit below still applies.
# ---------------------------------------------------------------------
try {
    $handoutMaster = $p.HandoutMaster
    $handoutMaster.HeadersFooters.DateAndTime.Value = "1/12/21"
} catch {
}

# ---------------------------------------------------------------------
# Slide 3 holds the two callouts that were copy-edited.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(3)

# ---------------------------------------------------------------------
# 2) Typo fix: "Consortium aggres on CC Def" -> "Consortium agrees on CC Def"
#    That sentence is its own run inside the "Rounded Rectangular
#    Callout 51" shape, so replace just that substring to keep the
#    run's own formatting (and the sibling "1. " run) untouched.
# ---------------------------------------------------------------------
$calloutCCDef = $slide.Shapes.Item("Rounded Rectangular Callout 51")
$rangeCCDef = $calloutCCDef.TextFrame.TextRange
$oldSentence = "Consortium aggres on CC Def"
$newSentence = "Consortium agrees on CC Def"
$matchStart = $rangeCCDef.Text.IndexOf($oldSentence)
if ($matchStart -ge 0) {
    $sentenceRange = $rangeCCDef.Characters($matchStart + 1, $oldSentence.Length)
    $sentenceRange.Text = $newSentence
}

# ---------------------------------------------------------------------
# 3) Renumber step 12 -> 13 in the "Rounded Rectangular Callout 80"
#    shape. The old text is a single run "12. Submit registerEnclave tx
#    for ordering"; the new text keeps "Submit registerEnclave tx for
#    ordering" as that run and puts "13. " in a new leading run.
# ---------------------------------------------------------------------
$calloutRegisterEnclave = $slide.Shapes.Item("Rounded Rectangular Callout 80")
$rangeRegisterEnclave = $calloutRegisterEnclave.TextFrame.TextRange
$oldPrefix = "12. "
if ($rangeRegisterEnclave.Text.IndexOf($oldPrefix) -eq 0) {
    $prefixRange = $rangeRegisterEnclave.Characters(1, $oldPrefix.Length)
    $prefixRange.Text = ""
    $rangeRegisterEnclave.InsertBefore("13. ") | Out-Null
}
